$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.769.51"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.437.53"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.97"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.70"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.478"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.63"
$ws.Range("E9").Value = "  +5.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.122"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.382"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.056.38"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.467.13"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000176"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.780.68"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.72"
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.84"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.63"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.27"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.48"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.560"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.584.74"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.03"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000110"
$ws.Range("E27").Value = "  -5.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").Value = "  -4.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.155"
$ws.Range("E31").Value = "  +3.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.87"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.41"
$ws.Range("E33").Value = "  -4.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.487.54"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.77"
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.15"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.70"
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.48"
$ws.Range("E40").Value = "  -4.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0767"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.790"
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.31"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.31"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.60"
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.27"
$ws.Range("E47").Value = "  -7.88%  "
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.65"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.310.95"
$ws.Range("E50").Value = "  -6.05%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.884"
$ws.Range("E51").Value = "  -0.93%  "
